$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$win = $excel.ActiveWindow
for ($i=1; $i -le $win.Panes.Count; $i++) {
  $p = $win.Panes.Item($i)
  Write-Host "Pane $i scrollrow=" $p.ScrollRow " scrollcol=" $p.ScrollColumn
}
